$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
    "23-9=14",
    "69+23=92",
    "91-79=12",
    "85-3=82",
    "5+54=59",
    "91-61=30",
    "1+43=44",
    "6-4=2",
    "34+51=85",
    "15+7=22",
    "41+21=62",
    "49+4=53",
    "87-29=58",
    "65-46=19",
    "94-56=38",
    "55+11=66",
    "79+11=90",
    "23+50=73",
    "52+31=83",
    "74+20=94",
    "79-26=53",
    "12-4=8",
    "89-69=20",
    "81+16=97",
    "99-21=78",
    "7+80=87",
    "45+35=80",
    "72-55=17",
    "99-89=10",
    "88-5=83",
    "37+25=62",
    "26+53=79",
    "12+82=94",
    "29-0=29",
    "65-22=43",
    "9+83=92",
    "64+25=89",
    "21+29=50",
    "91+7=98",
    "95-6=89",
    "27+33=60",
    "73+23=96",
    "64-57=7",
    "56-30=26",
    "72-38=34",
    "48+1=49",
    "87-2=85",
    "35-20=15",
    "9+70=79",
    "26-5=21",
    "19+54=73",
    "43-7=36",
    "90-90=0",
    "8+21=29",
    "60-44=16",
    "92-11=81",
    "30+35=65",
    "17+55=72",
    "62-59=3",
    "82-68=14",
    "90+5=95",
    "36+28=64",
    "78-73=5",
    "32+53=85",
    "85-17=68",
    "89-10=79",
    "14+78=92",
    "88-69=19",
    "22+18=40",
    "36+34=70",
    "7+92=99",
    "6+25=31",
    "12+74=86",
    "60-53=7",
    "80-4=76",
    "15+34=49",
    "14-6=8",
    "17+65=82",
    "75-22=53",
    "7+22=29",
    "29+63=92",
    "66+5=71",
    "98-77=21",
    "48+15=63",
    "17+6=23",
    "25+59=84",
    "52-9=43",
    "16+51=67",
    "79-41=38",
    "63+8=71",
    "78+6=84",
    "51+4=55",
    "81-32=49",
    "69-2=67",
    "42+50=92",
    "88-1=87",
    "28-4=24",
    "88-47=41",
    "68-5=63",
    "91-9=82"
)
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}
